$d = $word.ActiveDocument

$replacements = @(
    @('720×2=1440', '646×5=3230'),
    @('259×9=2331', '717×5=3585'),
    @('496×2=992', '683×5=3415'),
    @('755×6=4530', '291×8=2328'),
    @('779×4=3116', '947×4=3788'),
    @('702×7=4914', '653×2=1306'),
    @('378×2=756', '665×6=3990'),
    @('122×3=366', '272×9=2448'),
    @('898×6=5388', '620×5=3100'),
    @('191×7=1337', '130×7=910'),
    @('349×7=2443', '514×2=1028'),
    @('661×3=1983', '911×8=7288'),
    @('130×4=520', '482×8=3856'),
    @('411×6=2466', '438×9=3942'),
    @('799×4=3196', '844×6=5064'),
    @('446×6=2676', '307×2=614'),
    @('543×8=4344', '732×2=1464'),
    @('602×2=1204', '370×4=1480'),
    @('424×8=3392', '830×4=3320'),
    @('939×8=7512', '400×2=800'),
    @('943×2=1886', '616×5=3080'),
    @('650×7=4550', '725×8=5800'),
    @('578×3=1734', '339×2=678'),
    @('480×8=3840', '197×6=1182'),
    @('855×8=6840', '848×7=5936'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $old"
    }
}

$d.Save()
